$p = $ppt.ActivePresentation

$notes = @("Ben", "Ben", "Ross", "Ben", "Ross", "Ben", "Ross", "Ben", "Ross", "Ben", "Ben")

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $np = $s.NotesPage
    $body = $np.Shapes.Placeholders.Item(2)
    $body.TextFrame.TextRange.Text = $notes[$i - 1]
}
